# Refresh the "Elapsed Duration(Hrs)" column (G) on each region sheet.
# These are plain text HH:MM:SS-style durations (NOW() - PCM Created At,
# captured as a static snapshot rather than a live formula), so we just
# overwrite the stale strings with their recalculated values.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("R1").Range("G2").Value = "3926:29:43"
$wb.Worksheets.Item("R1").Range("G3").Value = "66:02:21"

$wb.Worksheets.Item("R2").Range("G2").Value = "12107:53:23"
$wb.Worksheets.Item("R2").Range("G3").Value = "3237:36:52"
$wb.Worksheets.Item("R2").Range("G4").Value = "475:48:26"

$wb.Worksheets.Item("R4").Range("G2").Value = "2953:43:12"
$wb.Worksheets.Item("R4").Range("G3").Value = "180:55:27"

$wb.Worksheets.Item("R5").Range("G2").Value = "427:42:11"

$wb.Worksheets.Item("R6").Range("G2").Value = "68:14:29"
